$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) interpretation for cells being updated, so that
# numeric-looking strings (e.g. "572.45", "0.0000158") are preserved as
# literal text, matching the original inlineStr cell contents.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.550.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.173.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.76%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.76%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -6.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.179.66"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.56%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.82"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.725.41"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.77%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.590.91"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000158"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.47%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.181.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "417.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.52%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.56%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.68"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.204"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.14%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000105"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.80"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.62%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.77"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.07"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.38"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.80%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.76%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.99"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.94%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.71"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.48%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.704.28"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.24"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.33%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.33"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.44%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.11"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.85%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.717"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.83%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.57"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.28%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.44%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "293.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.78%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.55"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.31%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -13.29%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.25%  "
